# "Generate Report for Handback" - populate the Latest Target File / Latest
# Handback File columns (F/G) for each language sheet, flip the Status text
# from "Ready for handoff" to the handed-back message, and stamp the
# handback datetime.

$wb = $excel.ActiveWorkbook

$mdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/c20393397f88d0e354d17e761fe91d121ee0b23f/e2e/c075ce48-2946-4b29-9cd8-dc9e82ee0cbd.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aeb473f4a738f8adc03766b1c6fc57189819da1c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c075ce48-2946-4b29-9cd8-dc9e82ee0cbd.e4e5a9398226a09db97282e15f3993711e8ca2fd.zh-cn.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9cc85c07138dac0337d01fac5588719210712e45/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c075ce48-2946-4b29-9cd8-dc9e82ee0cbd.e4e5a9398226a09db97282e15f3993711e8ca2fd.de-de.xlf"

$targetFileName = "c075ce48-2946-4b29-9cd8-dc9e82ee0cbd.md"
$zhHandbackFileName = "c075ce48-2946-4b29-9cd8-dc9e82ee0cbd.e4e5a9398226a09db97282e15f3993711e8ca2fd.zh-cn.xlf"
$deHandbackFileName = "c075ce48-2946-4b29-9cd8-dc9e82ee0cbd.e4e5a9398226a09db97282e15f3993711e8ca2fd.de-de.xlf"

$handedBackStatus = "Handed back: in sync with en-US"

# ---- zh-cn sheet -----------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $handedBackStatus
$zh.Range("C3").Value = $handedBackStatus

$zh.Hyperlinks.Add($zh.Range("F2"), $mdUrl, "", "", $targetFileName)
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhHandbackFileName)
$zh.Hyperlinks.Add($zh.Range("F3"), $mdUrl, "", "", $targetFileName)
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhHandbackFileName)

$zh.Range("H2").Value = "2016-03-20 22:56:01"
$zh.Range("H3").Value = "2016-03-20 22:56:01"

# ---- de-de sheet -------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $handedBackStatus
$de.Range("C3").Value = $handedBackStatus

$de.Hyperlinks.Add($de.Range("F2"), $mdUrl, "", "", $targetFileName)
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deHandbackFileName)
$de.Hyperlinks.Add($de.Range("F3"), $mdUrl, "", "", $targetFileName)
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deHandbackFileName)

$de.Range("H2").Value = "2016-03-20 22:56:07"
$de.Range("H3").Value = "2016-03-20 22:56:07"

Write-Host "Handback report generated"
